$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "汽車" (cars) sheet

# --- Header row (row 1) ---------------------------------------------------
# Existing header cells B1:G1 get new "field name" header labels (style s=1
# is already applied to them, so plain Value assignment keeps it).
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"

# New header cells H1:N1 -- copy formatting from B1 (style s=1) then set text.
$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,8))
$ws.Cells.Item(1,8).Value = "property_category"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,9))
$ws.Cells.Item(1,9).Value = "category"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,10))
$ws.Cells.Item(1,10).Value = "date"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,11))
$ws.Cells.Item(1,11).Value = "legislator_name"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,12))
$ws.Cells.Item(1,12).Value = "legislator_id"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,13))
$ws.Cells.Item(1,13).Value = "source_file"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,14))
$ws.Cells.Item(1,14).Value = "index"

# --- Data row (row 2) ------------------------------------------------------
# B2:G2 already hold the correct car data and keep their values unchanged.

# New data cells H2:N2 -- copy formatting from B2 (style s=2) then set values.
$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,8))
$ws.Cells.Item(2,8).Value = "land"

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,9))
$ws.Cells.Item(2,9).Value = "normal"

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,10))
$ws.Cells.Item(2,10).NumberFormat = "@"   # force text so the date-looking string is not auto-converted
$ws.Cells.Item(2,10).Value = "2013-12-17"

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,11))
$ws.Cells.Item(2,11).Value = "蔡煌瑯"

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,12))
$ws.Cells.Item(2,12).Value = 752

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,13))
$ws.Cells.Item(2,13).Value = "tmpc9fc1"

$ws.Cells.Item(2,2).Copy($ws.Cells.Item(2,14))
$ws.Cells.Item(2,14).Value = 35

Write-Host "sheet3 (cars) updated with capacity/common columns"
